$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("year")

# Row 2 currently stores the combined "2006, 2010" text in E2.
# Split it into two proper year rows: E2 -> 2006 (numeric), and a new
# row 3 duplicating the rest of row 2 with E3 -> 2010 (numeric).
$ws.Cells.Item(2, 5).Value = 2006

$ws.Cells.Item(3, 1).Value = "grc"
$ws.Cells.Item(3, 2).Value = "all"
$ws.Cells.Item(3, 3).Value = "NA"
$ws.Cells.Item(3, 4).Value = "NA"
$ws.Cells.Item(3, 5).Value = 2010
$ws.Cells.Item(3, 6).Value = "NA"
$ws.Cells.Item(3, 7).Value = "NA"

# Match the number formatting used on E2 (thousands-style numFmt) for E3.
$ws.Cells.Item(3, 5).NumberFormat = $ws.Cells.Item(2, 5).NumberFormat

[void]$ws.Range("E4").Select()
